# "add Use Item bug"
#
# NPC.xlsx / Property sheet:
#  - Row 49 (ShowName / 显示名字) loses its heavy custom row/cell formatting
#    (the row had been using the "highlight" style block s="6"/s="7"); it
#    becomes a plain data row like the rows around it, keeping only the
#    text-format ("@") on column I (the "Friend" scope column).
#  - Two brand-new attribute rows are appended after EquipIDRef:
#       Icon      / 图标        (row 51)
#       ShowCard  / 卡牌背景     (row 52)
#  - The sheet's view is scrolled down and the selection moved to A48.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-AttrRow($row, $idName, $desc) {
    $ws.Cells.Item($row, 1).Value = $idName
    $ws.Cells.Item($row, 2).Value = "string"
    $ws.Cells.Item($row, 3).Value = $false
    $ws.Cells.Item($row, 4).Value = $false
    $ws.Cells.Item($row, 5).Value = $false
    $ws.Cells.Item($row, 6).Value = $true
    $ws.Cells.Item($row, 7).Value = 0
    $ws.Cells.Item($row, 8).Value = 0
    $ws.Cells.Item($row, 9).Value = "Friend"
    $ws.Cells.Item($row, 9).NumberFormat = "@"
    $ws.Cells.Item($row, 10).Value = $desc
}

# Row 49 keeps its data (ShowName / 显示名字) but drops the special
# highlighted formatting that used to mark it - strip the row/cell styles
# first, then restore the single text-format cell (column I) it still needs.
$ws.Rows(49).ClearFormats()
Set-AttrRow 49 "ShowName" "显示名字"

# New rows for the item-card related attributes.
Set-AttrRow 51 "Icon" "图标"
Set-AttrRow 52 "ShowCard" "卡牌背景"

# Update the view: scroll further down and move the active selection.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 25
$win.ScrollColumn = 1
$ws.Range("A48").Select()
